$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = 5
    6  = -1
    9  = 0
    18 = 0
    19 = -1
    22 = -4
    23 = -8
    24 = -1
    25 = -2
    26 = 2
    29 = -6
    30 = -7
    32 = 7
    34 = 1
    35 = 1
    36 = -5
    37 = -2
    39 = -1
    40 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
